$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header/labels) - values in B1:E1 were changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (Subj 1 data) - values in B2:E2 were changed
$ws.Range("B2").Value = 11.088586800610617
$ws.Range("C2").Value = 10.975468096729349
$ws.Range("D2").Value = 12.009443207668815
$ws.Range("E2").Value = 11.247681488554406

# Row 3 (Subj 2 data) - values in B3:E3 were changed
$ws.Range("B3").Value = 10.319760442308223
$ws.Range("C3").Value = 9.0770161858068636
$ws.Range("D3").Value = 10.672956785928051
$ws.Range("E3").Value = 11.102607522646252

# Update the selected range shown when the sheet is reopened
$ws.Range("B1:E3").Select()
